# Update Daily Report: 2026-02-09
# Appends the new daily snapshot (date serial 46059 / 2026-02-06) for every
# depository/region-type combination to Daily_Data, then refreshes the
# derived Today_Summary and Monthly_Stats sheets to reflect the new data.

$wb = $excel.ActiveWorkbook

$daily = $wb.Worksheets.Item("Daily_Data")
$today = $wb.Worksheets.Item("Today_Summary")
$monthly = $wb.Worksheets.Item("Monthly_Stats")

# ---------------------------------------------------------------------------
# 1. Append the new date group (46059) to Daily_Data, rows 530-551.
#    Columns: A=Date, B=Region_Type, C=PREV_TOTAL, D=RECEIVED, E=WITHDRAWN,
#             F=NET_CHANGE, G=ADJUSTMENT, H=TOTAL_TODAY
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Row = 530; Name = "ASAHI DEPOSITORY LLC Registered";                              C = 0;          D = 0; E = 0;         F = 0;          G = 0;       H = 0 },
    @{ Row = 531; Name = "ASAHI DEPOSITORY LLC Eligible";                                C = 0;          D = 0; E = 0;         F = 0;          G = 0;       H = 0 },
    @{ Row = 532; Name = "BRINK'S, INC. Registered";                                     C = 76497.842;  D = 0; E = 0;         F = 0;          G = -874.54; H = 75623.302 },
    @{ Row = 533; Name = "BRINK'S, INC. Eligible";                                       C = 82678.788;  D = 0; E = 0;         F = 0;          G = 874.54;  H = 83553.32799999999 },
    @{ Row = 534; Name = "CNT DEPOSITORY, INC. Registered";                              C = 1246.06;    D = 0; E = 0;         F = 0;          G = 0;       H = 1246.06 },
    @{ Row = 535; Name = "CNT DEPOSITORY, INC. Eligible";                                C = 0;          D = 0; E = 0;         F = 0;          G = 0;       H = 0 },
    @{ Row = 536; Name = "DELAWARE DEPOSITORY Registered";                               C = 1633.941;   D = 0; E = 0;         F = 0;          G = 0;       H = 1633.941 },
    @{ Row = 537; Name = "DELAWARE DEPOSITORY Eligible";                                 C = 18459.584;  D = 0; E = 0;         F = 0;          G = 0;       H = 18459.584 },
    @{ Row = 538; Name = "HSBC BANK, USA Registered";                                    C = 1394.758;   D = 0; E = 0;         F = 0;          G = 0;       H = 1394.758 },
    @{ Row = 539; Name = "HSBC BANK, USA Eligible";                                      C = 9281.978999999999; D = 0; E = 0;   F = 0;          G = 0;       H = 9281.978999999999 },
    @{ Row = 540; Name = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered";     C = 2395.448;   D = 0; E = 0;         F = 0;          G = 0;       H = 2395.448 },
    @{ Row = 541; Name = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible";       C = 0;          D = 0; E = 0;         F = 0;          G = 0;       H = 0 },
    @{ Row = 542; Name = "JP MORGAN CHASE BANK NA Registered";                           C = 114985.579; D = 0; E = 0;         F = 0;          G = 0;       H = 114985.579 },
    @{ Row = 543; Name = "JP MORGAN CHASE BANK NA Eligible";                             C = 75484.511;  D = 0; E = 0;         F = 0;          G = 0;       H = 75484.511 },
    @{ Row = 544; Name = "LOOMIS INTERNATIONAL (US) LLC Registered";                     C = 63745.991;  D = 0; E = 0;         F = 0;          G = 0;       H = 63745.991 },
    @{ Row = 545; Name = "LOOMIS INTERNATIONAL (US) LLC Eligible";                       C = 132077.206; D = 0; E = 63071.566; F = -63071.566; G = 0;       H = 69005.64 },
    @{ Row = 546; Name = "MALCA-AMIT USA, LLC Registered";                               C = 395.145;    D = 0; E = 0;         F = 0;          G = 0;       H = 395.145 },
    @{ Row = 547; Name = "MALCA-AMIT USA, LLC Eligible";                                 C = 0;          D = 0; E = 0;         F = 0;          G = 0;       H = 0 },
    @{ Row = 548; Name = "MANFRA, TORDELLA & BROOKES, LLC Registered";                   C = 50220.42;   D = 0; E = 0;         F = 0;          G = 0;       H = 50220.42 },
    @{ Row = 549; Name = "MANFRA, TORDELLA & BROOKES, LLC Eligible";                      C = 1804.683;   D = 0; E = 0;         F = 0;          G = 0;       H = 1804.683 },
    @{ Row = 550; Name = "STONEX PRECIOUS METALS LLC Registered";                        C = 14122.765;  D = 0; E = 0;         F = 0;          G = 0;       H = 14122.765 },
    @{ Row = 551; Name = "STONEX PRECIOUS METALS LLC Eligible";                          C = 16.075;     D = 0; E = 0;         F = 0;          G = 0;       H = 16.075 }
)

$dateSerial = 46059

foreach ($item in $newRows) {
    $r = $item.Row
    $daily.Range("A$r").Value = $dateSerial
    $daily.Range("A$r").NumberFormat = $daily.Range("A$($r - 1)").NumberFormat
    $daily.Range("B$r").Value = $item.Name
    $daily.Range("C$r").Value = $item.C
    $daily.Range("D$r").Value = $item.D
    $daily.Range("E$r").Value = $item.E
    $daily.Range("F$r").Value = $item.F
    $daily.Range("G$r").Value = $item.G
    $daily.Range("H$r").Value = $item.H
}

# ---------------------------------------------------------------------------
# 2. Refresh Today_Summary with the latest snapshot (46059) for the two
#    depositories whose figures changed: BRINK'S, INC. (row 3) and
#    LOOMIS INTERNATIONAL (US) LLC (row 9).
# ---------------------------------------------------------------------------

$today.Range("B3").Value = 83553.32799999999   # BRINK'S, INC. Eligible
$today.Range("C3").Value = 75623.302            # BRINK'S, INC. Registered

$today.Range("B9").Value = 69005.64             # LOOMIS INTERNATIONAL (US) LLC Eligible
$today.Range("D9").Value = 132751.631           # LOOMIS INTERNATIONAL (US) LLC Total_Stock

# ---------------------------------------------------------------------------
# 3. Refresh Monthly_Stats for February 2026.
# ---------------------------------------------------------------------------

# Month summary row (row 2, YearMonth = 2026-02)
$monthly.Range("B2").Value = 257605.8
$monthly.Range("C2").Value = 325763.409
$monthly.Range("D2").Value = 583369.209

# Detail rows (2026-02 section): RECEIVED/WITHDRAWN/TOTAL_TODAY per
# depository/region-type, updated for the companies affected by the new day.
$monthly.Range("E10").Value = 83553.32799999999   # BRINK'S, INC. Eligible - TOTAL_TODAY
$monthly.Range("E11").Value = 75623.302            # BRINK'S, INC. Registered - TOTAL_TODAY

$monthly.Range("D22").Value = 63071.566            # LOOMIS INTERNATIONAL (US) LLC Eligible - WITHDRAWN
$monthly.Range("E22").Value = 69005.64             # LOOMIS INTERNATIONAL (US) LLC Eligible - TOTAL_TODAY
